$d = $word.ActiveDocument

# Locate the "Філософія Відродження" heading paragraph (style "Заголовок-2" / -20).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Філософія Відродження" + [char]13) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Філософія Відродження' heading paragraph"
}

# Promote it to the top-level heading style and upper-case its text.
$target.Style = "Заголовок-1"
$target.Range.Text = "ФІЛОСОФІЯ ВІДРОДЖЕННЯ"

# Re-acquire the (now retitled) heading paragraph and find the paragraph right after it
# (the pre-existing empty paragraph that follows the heading).
$afterHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ФІЛОСОФІЯ ВІДРОДЖЕННЯ" + [char]13) {
        $afterHeading = $p.Next()
        break
    }
}

if ($afterHeading -eq $null) {
    throw "Could not locate the paragraph following the heading"
}

$insertionPoint = $afterHeading.Range
$insertionPoint.Collapse(1)

$xml = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>1. Філософія Відродження:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Виникла в 14 столітті в Італії і поширилася по всій Європі.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Відкинула середньовічну схоластику і зосередилася на гуманізмі, індивідуалізмі та світськості.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Гуманізм підкреслював цінність і потенціал кожної людини та її досягнень.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Індивідуалізм наголошував на індивідуальній свободі, творчості та самовираженні.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Секуляризм</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> підкреслював важливість мирських справ над релігійними.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- Серед відомих філософів епохи Відродження - </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Нікколо</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Макіавеллі, Френсіс Бекон і Мішель де Монтень.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>2. Філософія Нового часу:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Виникла в 17 столітті як відповідь на епоху Відродження та наукову революцію.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Наголошувала на розумі, емпіричному спостереженні та скептицизмі.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- До відомих філософів Нового часу належать </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Рене</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Декарт, Джон Локк та </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Іммануїл</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Кант.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Філософія Декарта наголошує на важливості сумніву та силі розуму для отримання певного знання.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Філософія Локка наголошувала на важливості досвіду і відкидала вроджені ідеї.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>- Філософія Канта наголошувала на обмеженості розуму та важливості морального обов'язку.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>3. Філософія Просвітництва:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Виникла у 18 столітті і наголошувала на розумі, прогресі та свободі.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Критикувала традиційну владу і виступала за індивідуальну свободу і соціальну рівність.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Серед відомих філософів Просвітництва - Жан-Жак Руссо, Вольтер і Джон Стюарт Мілль.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Філософія Руссо підкреслювала важливість суспільного договору та загальної волі народу.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Філософія Вольтера наголошувала на важливості толерантності, свободи слова та відокремлення церкви від держави.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t>- Філософія Мілля наголошувала на важливості індивідуальної свободи та неприйнятті тиранії.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="851"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Чи є </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>протирічними</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> теорії </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Р.Декарта</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> і </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ф.Бекона</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> щодо пізнавального процесу у людини, та </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Б.Спінози</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> та Г. Лейбніца щодо існування світу? </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Обгрунтуйте</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> свою відповідь</w:t>
      </w:r>
    </w:p><w:p><w:pPr><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:ind w:firstLine="851"/></w:pPr></w:p></pkg:xmlData>
'@

[void]$insertionPoint.InsertXML($xml)

Write-Output "OK"
